$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-style the gem cells (H2:H4) from the "Bad" style (red) to the blue
# style already used elsewhere in the sheet (same style as A2/C2/D2/I2),
# then reorder the gem names shown in that column.
$ws.Range("A2").Copy()
$ws.Range("H2:H4").PasteSpecial(-4122)

$ws.Range("H2").Value = "saphire"
$ws.Range("H3").Value = "emerald"
$ws.Range("H4").Value = "ruby"

# --- Add the two new entries (new shared strings get appended in the order
# they are first used, so "potion" must be written before "pendulum" to land
# on shared-string indices 70/71 respectively).
$ws.Range("A5").Value = "potion"
$ws.Range("C4").Value = "pendulum"

# Give the two new cells a brand-new fill (a green tint, distinct from the
# other colored columns). We first copy an existing solid-filled cell's
# format over so the Interior already has a solid pattern + explicit color;
# only then do we swap the theme color, which keeps the engine from
# generating a spurious intermediate "blank" fill entry.
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("C4").PasteSpecial(-4122)

$ws.Range("A5").Interior.ThemeColor = 10
$ws.Range("C4").Interior.ThemeColor = 10

# restore the values (PasteSpecial of formats only shouldn't have touched
# them, but make sure they are what we expect)
$ws.Range("A5").Value = "potion"
$ws.Range("C4").Value = "pendulum"

# --- Update the active selection shown when the sheet is reopened.
[void]$ws.Range("D6").Select()

Write-Host "edit complete"
